$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: B1 stays "wastage" (string table will shift since "unsure" gets pruned);
# add new C1 header "enhanced wastage"
$ws.Range("C1").Value = "enhanced wastage"

# Populate column C ("enhanced wastage") for each data row
$ws.Cells.Item(2, 3).Value = 0.05
$ws.Cells.Item(3, 3).Value = 0.05
$ws.Cells.Item(4, 3).Value = 0.05
$ws.Cells.Item(5, 3).Value = 0.015
$ws.Cells.Item(6, 3).Value = 0.015
$ws.Cells.Item(7, 3).Value = 0.025
$ws.Cells.Item(8, 3).Value = 0.05
$ws.Cells.Item(9, 3).Value = 0.025
$ws.Cells.Item(10, 3).Value = 0.05
$ws.Cells.Item(11, 3).Value = 0.05
$ws.Cells.Item(12, 3).Value = 0.025
$ws.Cells.Item(13, 3).Value = 0.015
$ws.Cells.Item(14, 3).Value = 0.05
$ws.Cells.Item(15, 3).Value = 0.015
$ws.Cells.Item(16, 3).Value = 0.075
$ws.Cells.Item(17, 3).Value = 0.05
$ws.Cells.Item(18, 3).Value = 0.005
$ws.Cells.Item(19, 3).Value = 0.025
$ws.Cells.Item(20, 3).Value = 0.075
$ws.Cells.Item(21, 3).Value = 0.05
$ws.Cells.Item(22, 3).Value = 0.05
$ws.Cells.Item(23, 3).Value = 0.05
$ws.Cells.Item(24, 3).Value = 0.075
$ws.Cells.Item(25, 3).Value = 0.05
$ws.Cells.Item(26, 3).Value = 0.075
$ws.Cells.Item(27, 3).Value = 0.025
$ws.Cells.Item(28, 3).Value = 0.05
$ws.Cells.Item(29, 3).Value = 0.025
$ws.Cells.Item(30, 3).Value = 0.075
$ws.Cells.Item(31, 3).Value = 0.075
$ws.Cells.Item(32, 3).Value = 0.075
$ws.Cells.Item(33, 3).Value = 0.05
$ws.Cells.Item(34, 3).Value = 0.05
$ws.Cells.Item(35, 3).Value = 0.025
$ws.Cells.Item(36, 3).Value = 0.05
$ws.Cells.Item(37, 3).Value = 0.05
$ws.Cells.Item(38, 3).Value = 0.025
$ws.Cells.Item(39, 3).Value = 0.05
$ws.Cells.Item(40, 3).Value = 0.05
$ws.Cells.Item(41, 3).Value = 0.05
$ws.Cells.Item(42, 3).Value = 0.05
$ws.Cells.Item(43, 3).Value = 0.05
$ws.Cells.Item(44, 3).Value = 0.05
$ws.Cells.Item(45, 3).Value = 0.025
$ws.Cells.Item(46, 3).Value = 0.05
$ws.Cells.Item(47, 3).Value = 0.005
$ws.Cells.Item(48, 3).Value = 0.005
$ws.Cells.Item(49, 3).Value = 0.075
$ws.Cells.Item(50, 3).Value = 0.075

# The three rows that previously held the placeholder "unsure" text in column C
# (rows 19, 29, 40) are overwritten above with numeric enhanced-wastage values,
# which makes the "unsure" shared string unused and it is dropped on save.

# Update the view: clear the scrolled topLeftCell and move the selection to D6
[void]$ws.Range("D6").Select()
